$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12; existing rows 12..27 shift down to 13..28,
# preserving their data/formatting.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with this week's record.
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 44650
$ws.Cells.Item(12, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112045
$ws.Cells.Item(12, 7).Value = "Zapallo"
$ws.Cells.Item(12, 8).Value = "Camote"
$ws.Cells.Item(12, 9).Value = "2a (cosecha)"
$ws.Cells.Item(12, 10).Value = 1300
$ws.Cells.Item(12, 11).Value = 400
$ws.Cells.Item(12, 12).Value = 430
$ws.Cells.Item(12, 13).Value = 415
$ws.Cells.Item(12, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(12, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 16).Value = 415
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = "Hortaliza"
